# US 3.3 commit files
# - "About" sheet: remove the "Notes from Ben:" block (rows 12-15)
# - "BDPbES" sheet:
#     * add a header label in A1 ("Priority Order (dimensionless)"),
#       bold + wrap text, with a taller row 1
#     * rename the "coal to gas" row (A13) to "lignite"
#     * append three new fuel-type rows (crude oil, heavy or residual fuel
#       oil, municipal solid waste) with the same priority (2) and the same
#       $B<row> fill-across formula pattern used by every other row
# - make "About" the active/selected sheet (it was "BDPbES" before)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# About sheet: drop the "Notes from Ben" discussion (rows 12-15)
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Rows("12:15").Delete()

# ---------------------------------------------------------------------------
# BDPbES sheet: rename, new rows, then header cell
# ---------------------------------------------------------------------------
$bd = $wb.Worksheets.Item("BDPbES")

# "coal to gas" -> "lignite"
$bd.Range("A13").Value = "lignite"

# New fuel rows 15-17, same pattern as the existing rows: B = priority,
# C:AK = "=$B<row>" fill across all the year columns
$newRows = @(
    @{ Row = 15; Name = "crude oil" },
    @{ Row = 16; Name = "heavy or residual fuel oil" },
    @{ Row = 17; Name = "municipal solid waste" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $bd.Range("A$r").Value = $item.Name
    $bd.Range("B$r").Value = 2
    $bd.Range("C$r`:AK$r").Formula = "=`$B$r"
}

# New header above the year columns
$bd.Range("A1").Value = "Priority Order (dimensionless)"
$bd.Range("A1").Font.Bold = $true
$bd.Range("A1").WrapText = $true
$bd.Rows(1).RowHeight = 30

# ---------------------------------------------------------------------------
# "About" becomes the selected/active sheet (was "BDPbES")
# ---------------------------------------------------------------------------
$about.Activate()
